# Auto commit at 2025-11-03  8:34:06.61
#
# Refresh the daily metrics snapshot:
#  - "Metrics" sheet: new period totals (B2:B13)
#  - "today" sheet: today's incremental figures (B3:B6) and the
#    "cumulative + today" formulas in F11:F14 (filled down, relative refs)
#  - view state: leave "Metrics" as the active sheet / tab, with the
#    "today" sheet's last touched range selected, matching the author's
#    on-screen state when they saved.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Metrics sheet: refreshed period totals
# ---------------------------------------------------------------------
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 25539.82
$metrics.Range("B3").Value  = 22365.42
$metrics.Range("B4").Value  = 7764.92
$metrics.Range("B5").Value  = 1068
$metrics.Range("B6").Value  = 4821785.5699999994
$metrics.Range("B7").Value  = 4064442.1000000006
$metrics.Range("B8").Value  = 1414724.75
$metrics.Range("B9").Value  = 187275
$metrics.Range("B10").Value = 33287166.560000002
$metrics.Range("B11").Value = 31339717.259999998
$metrics.Range("B12").Value = 11696446.790000001
$metrics.Range("B13").Value = 1284905

# ---------------------------------------------------------------------
# today sheet: fill in today's figures + "cumulative+today" formulas
# ---------------------------------------------------------------------
$today = $wb.Worksheets.Item("today")

$today.Range("B3").Value = 13177.81
$today.Range("B4").Value = 11622.74
$today.Range("B5").Value = 3986.97
$today.Range("B6").Value = 528

$today.Range("F11").Formula = "=E11+B3"
$today.Range("F12:F14").Formula = "=E12+B4"

# ---------------------------------------------------------------------
# View state: select the recomputed range on "today" first, then make
# "Metrics" the active sheet/tab with D7 selected (matches the state the
# workbook was saved in).
# ---------------------------------------------------------------------
[void]$today.Range("F11:F22").Select()

[void]$metrics.Activate()
[void]$metrics.Range("D7").Select()
